$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5529
$ws.Range("L3").Value = 6017
$ws.Range("L4").Value = 1478
$ws.Range("L5").Value = 360
$ws.Range("L6").Value = 4939
$ws.Range("L7").Value = 18323

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 593
$ws.Range("L8").Value = 1213
$ws.Range("L10").Value = 123
$ws.Range("L11").Value = 299
$ws.Range("L12").Value = 42
$ws.Range("L15").Value = 147
$ws.Range("L16").Value = 39
$ws.Range("L19").Value = 498
$ws.Range("L21").Value = 58
$ws.Range("L23").Value = 199
$ws.Range("L25").Value = 112
$ws.Range("L27").Value = 160
$ws.Range("L29").Value = 1033
$ws.Range("L31").Value = 180
$ws.Range("L34").Value = 106
$ws.Range("L36").Value = 234
$ws.Range("L37").Value = 699
$ws.Range("L42").Value = 596
$ws.Range("L44").Value = 123
$ws.Range("L45").Value = 34
$ws.Range("L48").Value = 239
$ws.Range("L52").Value = 375
$ws.Range("L54").Value = 400
$ws.Range("L55").Value = 193
$ws.Range("L59").Value = 32
$ws.Range("L63").Value = 58
$ws.Range("L65").Value = 352
$ws.Range("L67").Value = 632
$ws.Range("L73").Value = 146
$ws.Range("L77").Value = 121
$ws.Range("L78").Value = 232
$ws.Range("L79").Value = 499
$ws.Range("L83").Value = 406
$ws.Range("L84").Value = 177
$ws.Range("L85").Value = 914
$ws.Range("L86").Value = 125
$ws.Range("L91").Value = 241
$ws.Range("L93").Value = 90
$ws.Range("L94").Value = 223
$ws.Range("L95").Value = 258
$ws.Range("L99").Value = 318
$ws.Range("L101").Value = 18323

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 203
$ws.Range("L7").Value = 593

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 115
$ws.Range("L3").Value = 90
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 277
$ws.Range("L3").Value = 375
$ws.Range("L4").Value = 53
$ws.Range("L7").Value = 914

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L5").Value = 9
$ws.Range("L7").Value = 375

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 354
$ws.Range("L3").Value = 427
$ws.Range("L4").Value = 88
$ws.Range("L7").Value = 1213

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 406

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 83
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 209
$ws.Range("L3").Value = 243
$ws.Range("L4").Value = 38
$ws.Range("L6").Value = 189
$ws.Range("L7").Value = 699

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 112
$ws.Range("L7").Value = 352

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 129
$ws.Range("L7").Value = 318

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 75
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 244
$ws.Range("L4").Value = 41
$ws.Range("L7").Value = 632

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 72
$ws.Range("L3").Value = 102
$ws.Range("L4").Value = 33
$ws.Range("L7").Value = 400

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 306
$ws.Range("L3").Value = 399
$ws.Range("L7").Value = 1033

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 239

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 154
$ws.Range("L6").Value = 139
$ws.Range("L7").Value = 498

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L4").Value = 51
$ws.Range("L7").Value = 596

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 232

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 193

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 112
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 162
$ws.Range("L7").Value = 499

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 83
$ws.Range("L7").Value = 234

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L3").Value = 25
$ws.Range("L4").Value = 67
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 39
